$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.435.34'
$ws.Range('E2').Value = '  +0.83%  '
$ws.Range('D3').Value = '3.520.75'
$ws.Range('E3').Value = '  +0.30%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '604.85'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('D6').Value = '174.08'
$ws.Range('E6').Value = '  +1.65%  '
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').Value = '3.516.34'
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('D10').Value = '0.197'
$ws.Range('E10').Value = '  -0.65%  '
$ws.Range('D11').Value = '7.28'
$ws.Range('E11').Value = '  +8.79%  '
$ws.Range('D12').Value = '0.589'
$ws.Range('E12').Value = '  +1.30%  '
$ws.Range('D13').Value = '46.25'
$ws.Range('E13').Value = '  -2.02%  '
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('D15').Value = '4.092.75'
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('B16').Value = 'BitcoinCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D16').Value = '615.63'
$ws.Range('E16').Value = '  -0.63%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').Value = '8.31'
$ws.Range('E17').Value = '  -0.62%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '70.442.59'
$ws.Range('E18').Value = '  +0.88%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.506.79'
$ws.Range('E19').Value = '  -0.18%  '
$ws.Range('E20').Value = '  +1.07%  '
$ws.Range('D21').Value = '17.57'
$ws.Range('E21').Value = '  +1.60%  '
$ws.Range('D22').Value = '0.883'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').Value = '9.07'
$ws.Range('E23').Value = '  -8.66%  '
$ws.Range('D24').Value = '99.54'
$ws.Range('E24').Value = '  +3.83%  '
$ws.Range('D25').Value = '15.67'
$ws.Range('E25').Value = '  -0.42%  '
$ws.Range('E26').Value = '  -2.87%  '
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').Value = '2.57'
$ws.Range('E28').Value = '  -0.76%  '
$ws.Range('D29').Value = '34.48'
$ws.Range('E29').Value = '  +4.03%  '
$ws.Range('D30').Value = '9.07'
$ws.Range('E30').Value = '  -1.42%  '
$ws.Range('E31').Value = '  -3.72%  '
$ws.Range('E32').Value = '  -2.63%  '
$ws.Range('D33').Value = '645.58'
$ws.Range('E33').Value = '  +13.75%  '
$ws.Range('E34').Value = '  -4.32%  '
$ws.Range('D35').Value = '6.84'
$ws.Range('E35').Value = '  -1.54%  '
$ws.Range('D36').Value = '0.0999'
$ws.Range('E36').Value = '  -1.11%  '
$ws.Range('B37').Value = 'Cosmos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D37').Value = '10.77'
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').Value = '3.54'
$ws.Range('E38').Value = '  +2.62%  '
$ws.Range('E39').Value = '  +6.59%  '
$ws.Range('D40').Value = '56.85'
$ws.Range('E40').Value = '  -0.37%  '
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range('E42').Value = '  +1.28%  '
$ws.Range('D43').Value = '0.0₃0743'
$ws.Range('E43').Value = '  +5.40%  '
$ws.Range('D44').Value = '3.367.46'
$ws.Range('E44').Value = '  +1.09%  '
$ws.Range('D45').Value = '0.310'
$ws.Range('E45').Value = '  -4.73%  '
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').Value = '2.91'
$ws.Range('E46').Value = '  -3.29%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '32.18'
$ws.Range('E47').Value = '  -2.45%  '
$ws.Range('E48').Value = '  -2.03%  '
$ws.Range('E49').Value = '  +1.17%  '
$ws.Range('D50').Value = '133.28'
$ws.Range('E50').Value = '  -1.42%  '
$ws.Range('E51').Value = '  -0.02%  '
